# Auto-generated Excel COM-interop edit script
# Applies cached-value updates to the Leve-profit sheets per the target diff.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 1645.8572
$ws.Range("J17").Value = 1645.8572
$ws.Range("L17").Value = 4937.571599999999
$ws.Range("N17").Value = -5273.571599999999
$ws.Range("H28").Value = 0
$ws.Range("I28").Value = 0
$ws.Range("J28").Value = 0
$ws.Range("K28").Value = 0
$ws.Range("L28").Value = 0
$ws.Range("M28").ClearContents()
$ws.Range("N28").ClearContents()
$ws.Range("H74").Value = 3960
$ws.Range("I74").Value = 3950
$ws.Range("K74").Value = 3950
$ws.Range("M74").Value = -3014
$ws.Range("H77").Value = 3960
$ws.Range("I77").Value = 3950
$ws.Range("K77").Value = 19750
$ws.Range("M77").Value = -15070
$ws.Range("H80").Value = 855.8889
$ws.Range("J80").Value = 837.875
$ws.Range("L80").Value = 2513.625
$ws.Range("N80").Value = -4509.625
$ws.Range("H83").Value = 855.8889
$ws.Range("J83").Value = 837.875
$ws.Range("L83").Value = 7540.875
$ws.Range("N83").Value = -17524.875
$ws.Range("H92").Value = 1194.1428
$ws.Range("I92").Value = 989.61536
$ws.Range("J92").Value = 3853
$ws.Range("K92").Value = 989.61536
$ws.Range("L92").Value = 3853
$ws.Range("M92").Value = 258.38464
$ws.Range("N92").Value = -6349
$ws.Range("H94").Value = 0
$ws.Range("I94").Value = 0
$ws.Range("K94").Value = 0
$ws.Range("M94").ClearContents()
$ws.Range("H100").Value = 0
$ws.Range("I100").Value = 0
$ws.Range("J100").Value = 0
$ws.Range("K100").Value = 0
$ws.Range("L100").Value = 0
$ws.Range("M100").ClearContents()
$ws.Range("N100").ClearContents()
$ws.Range("H107").Value = 0
$ws.Range("I107").Value = 0
$ws.Range("J107").Value = 0
$ws.Range("K107").Value = 0
$ws.Range("L107").Value = 0
$ws.Range("M107").ClearContents()
$ws.Range("N107").ClearContents()
$ws.Range("H118").Value = 1033.3334
$ws.Range("I118").Value = 1033.3334
$ws.Range("K118").Value = 3100.0002
$ws.Range("M118").Value = -1443.0002

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 10295.8
$ws.Range("I32").Value = 7036.4165
$ws.Range("J32").Value = 23333.334
$ws.Range("K32").Value = 7036.4165
$ws.Range("L32").Value = 23333.334
$ws.Range("M32").Value = -6749.4165
$ws.Range("N32").Value = -23907.334
$ws.Range("H74").Value = 3274.625
$ws.Range("I74").Value = 2671
$ws.Range("K74").Value = 2671
$ws.Range("M74").Value = -1797
$ws.Range("H77").Value = 3274.625
$ws.Range("I77").Value = 2671
$ws.Range("K77").Value = 13355
$ws.Range("M77").Value = -8987
$ws.Range("H93").Value = 0
$ws.Range("I93").Value = 0
$ws.Range("K93").Value = 0
$ws.Range("M93").ClearContents()
$ws.Range("H122").Value = 505
$ws.Range("I122").Value = 505
$ws.Range("K122").Value = 1515
$ws.Range("M122").Value = 935
$ws.Range("H124").Value = 20000
$ws.Range("J124").Value = 20000
$ws.Range("L124").Value = 20000
$ws.Range("N124").Value = -29820

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H64").Value = 1874
$ws.Range("I64").Value = 2000
$ws.Range("J64").Value = 1811
$ws.Range("K64").Value = 2000
$ws.Range("L64").Value = 1811
$ws.Range("M64").Value = -1775
$ws.Range("N64").Value = -2261
$ws.Range("H67").Value = 1874
$ws.Range("I67").Value = 2000
$ws.Range("J67").Value = 1811
$ws.Range("K67").Value = 2000
$ws.Range("L67").Value = 1811
$ws.Range("M67").Value = -1220
$ws.Range("N67").Value = -3371
$ws.Range("H80").Value = 797.875
$ws.Range("J80").Value = 815.3333
$ws.Range("L80").Value = 815.3333
$ws.Range("N80").Value = -2811.3333
$ws.Range("H83").Value = 797.875
$ws.Range("J83").Value = 815.3333
$ws.Range("L83").Value = 4076.6665
$ws.Range("N83").Value = -14060.6665
$ws.Range("H105").Value = 2272.7273
$ws.Range("I105").Value = 0
$ws.Range("K105").Value = 0
$ws.Range("M105").ClearContents()
$ws.Range("H107").Value = 2179.875
$ws.Range("I107").Value = 2156.5
$ws.Range("J107").Value = 2250
$ws.Range("K107").Value = 2156.5
$ws.Range("L107").Value = 2250
$ws.Range("M107").Value = -236.5
$ws.Range("N107").Value = -6090

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 3703
$ws.Range("I31").Value = 2447.5
$ws.Range("J31").Value = 4958.5
$ws.Range("K31").Value = 2447.5
$ws.Range("L31").Value = 4958.5
$ws.Range("M31").Value = -2152.5
$ws.Range("N31").Value = -5548.5
$ws.Range("H34").Value = 3703
$ws.Range("I34").Value = 2447.5
$ws.Range("J34").Value = 4958.5
$ws.Range("K34").Value = 2447.5
$ws.Range("L34").Value = 4958.5
$ws.Range("M34").Value = -2245.5
$ws.Range("N34").Value = -5362.5
$ws.Range("H99").Value = 5092.25
$ws.Range("I99").Value = 5092.25
$ws.Range("K99").Value = 5092.25
$ws.Range("M99").Value = -3594.25
$ws.Range("H103").Value = 13577.25
$ws.Range("I103").Value = 13577.25
$ws.Range("K103").Value = 13577.25
$ws.Range("M103").Value = -12405.25
$ws.Range("H124").Value = 149000
$ws.Range("J124").Value = 149000
$ws.Range("L124").Value = 149000
$ws.Range("N124").Value = -153910
$ws.Range("H126").Value = 5092.25
$ws.Range("I126").Value = 5092.25
$ws.Range("K126").Value = 15276.75
$ws.Range("M126").Value = -12806.75

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H3").Value = 5000
$ws.Range("I3").Value = 5000
$ws.Range("K3").Value = 15000
$ws.Range("M3").Value = -14888
$ws.Range("H6").Value = 27.142857
$ws.Range("I6").Value = 23.333334
$ws.Range("K6").Value = 70.00000199999999
$ws.Range("M6").Value = 42.99999800000001
$ws.Range("H107").Value = 1221
$ws.Range("I107").Value = 981
$ws.Range("K107").Value = 2943
$ws.Range("M107").Value = -1023
$ws.Range("H131").Value = 2023.2354
$ws.Range("I131").Value = 1500
$ws.Range("J131").Value = 2055.9375
$ws.Range("K131").Value = 4500
$ws.Range("L131").Value = 6167.8125
$ws.Range("M131").Value = 540
$ws.Range("N131").Value = -16247.8125
$ws.Range("H132").Value = 1449.5
$ws.Range("I132").Value = 900
$ws.Range("K132").Value = 8100
$ws.Range("M132").Value = -5570
$ws.Range("H134").Value = 466.66666
$ws.Range("I134").Value = 466.66666
$ws.Range("K134").Value = 1399.99998
$ws.Range("M134").Value = 3670.00002
$ws.Range("H140").Value = 461.91666
$ws.Range("I140").Value = 461.91666
$ws.Range("K140").Value = 1385.74998
$ws.Range("M140").Value = 3794.25002

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 29.333334
$ws.Range("I2").Value = 23.714285
$ws.Range("K2").Value = 23.714285
$ws.Range("M2").Value = 89.285715
$ws.Range("H27").Value = 4900
$ws.Range("J27").Value = 4900
$ws.Range("L27").Value = 4900
$ws.Range("N27").Value = -5232
$ws.Range("H80").Value = 22394.2
$ws.Range("J80").Value = 24982.834
$ws.Range("L80").Value = 24982.834
$ws.Range("N80").Value = -26978.834
$ws.Range("H83").Value = 22394.2
$ws.Range("J83").Value = 24982.834
$ws.Range("L83").Value = 124914.17
$ws.Range("N83").Value = -134898.17
$ws.Range("H102").Value = 264.16666
$ws.Range("I102").Value = 264.16666
$ws.Range("J102").Value = 0
$ws.Range("K102").Value = 264.16666
$ws.Range("L102").Value = 0
$ws.Range("M102").Value = 1357.83334
$ws.Range("N102").ClearContents()
$ws.Range("H113").Value = 1550
$ws.Range("I113").Value = 900
$ws.Range("K113").Value = 900
$ws.Range("M113").Value = 1270
$ws.Range("H122").Value = 2000
$ws.Range("I122").Value = 2000
$ws.Range("K122").Value = 6000
$ws.Range("M122").Value = -3550

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 3003
$ws.Range("I16").Value = 4516.2
$ws.Range("K16").Value = 4516.2
$ws.Range("M16").Value = -4346.2
$ws.Range("H40").Value = 0
$ws.Range("I40").Value = 0
$ws.Range("K40").Value = 0
$ws.Range("M40").ClearContents()
$ws.Range("H50").Value = 50001
$ws.Range("J50").Value = 50001
$ws.Range("L50").Value = 50001
$ws.Range("N50").Value = -51275
$ws.Range("H55").Value = 975.8889
$ws.Range("J55").Value = 895
$ws.Range("L55").Value = 895
$ws.Range("N55").Value = -1241

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H126").Value = 993.3333
$ws.Range("I126").Value = 993.3333
$ws.Range("K126").Value = 2979.9999
$ws.Range("M126").Value = -509.9998999999998
$ws.Range("H132").Value = 1251.9231
$ws.Range("I132").Value = 1106.909
$ws.Range("J132").Value = 2049.5
$ws.Range("K132").Value = 3320.727
$ws.Range("L132").Value = 6148.5
$ws.Range("M132").Value = -790.7270000000003
$ws.Range("N132").Value = -11208.5
$ws.Range("H136").Value = 1430.2222
$ws.Range("I136").Value = 1416.1666
$ws.Range("J136").Value = 1458.3334
$ws.Range("K136").Value = 4248.4998
$ws.Range("L136").Value = 4375.0002
$ws.Range("M136").Value = -1698.4998
$ws.Range("N136").Value = -9475.0002
$ws.Range("H138").Value = 100000
$ws.Range("J138").Value = 100000
$ws.Range("L138").Value = 100000
$ws.Range("N138").Value = -110280
